# "Got data import working; started on PDF generation"
#
# Target sheet is "|Toets 1" (already the active/selected sheet in the
# workbook), where a handful of imported data values changed and the
# cursor ended up on a different cell. A leftover conditional-format
# "dxf" (fill-only, from a Highlight-Cells-Rule that was created and then
# removed again) also shows up in the shared style table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet   # "|Toets 1" is already tabSelected/active

# --- Re-imported data values (row 8 and row 9) -----------------------
$ws.Range("F8").Value = 4
$ws.Range("G8").Value = 4
$ws.Range("H8").Value = 2
$ws.Range("I8").Value = 4
$ws.Range("I9").Value = 4

# --- Cursor moved to J10 before the file was saved --------------------
[void]$ws.Range("J10").Select()

# --- Orphaned conditional-format dxf (fill, theme accent6) ------------
# A "Highlight Cells Rule" with a themed fill was added and removed,
# which leaves its differential format behind in the style table even
# though no conditionalFormatting rule references it anymore.
$fc = $ws.Range("A1:I10").FormatConditions.Add(1, 3, "5")
$fc.Interior.Color = 3057486
$fc.Delete()
